# Weekly update: two new "Espárragos" (Vega Modelo de Temuco) price rows
# for the latest reporting week are inserted at the top of the existing
# data block (row 21), pushing the previously-reported rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 21 (existing rows 21:31
# shift down to 23:33).
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# --- New row 21 : "Primera" quality ---
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 44488
$ws.Range("D21").NumberFormat = $ws.Range("D23").NumberFormat
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 300000000
$ws.Range("G21").Value = "Espárragos"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 295
$ws.Range("K21").Value = 1300
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1415
$ws.Range("N21").Value = "$/kilo"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1415
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"

# --- New row 22 : "Segunda" quality ---
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44488
$ws.Range("D22").NumberFormat = $ws.Range("D23").NumberFormat
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 300000000
$ws.Range("G22").Value = "Espárragos"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = 1200
$ws.Range("N22").Value = "$/kilo"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 1200
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"
